$d = $word.ActiveDocument

# Remove the redundant "αστερισμό του" before "Αστερισμός μπότες" (Boötes),
# leaving the surrounding spaces intact (results in a double space). The
# search string includes "Αστερισμός" so the similar phrase referring to
# "αστερισμό του Ορίωνα" (Orion) elsewhere in the document is left untouched.
$d.Content.Find.Execute("αστερισμό του Αστερισμός", $true, $false, $false, $false, $false,
                         $true, 1, $false, " Αστερισμός", 2)
